# Timeline.xlsx update:
#  - Finished Face Recognition Engine items on the Timeline sheet
#    (fill in Actual/Finished dates + Additional Notes for several rows)
#  - Add a new "Development Module Breakdown" sheet tracking initial
#    development on the form app client (per-module start/end + unit test
#    status)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$dateFmt = "d-mmm-yy"

# ---------------------------------------------------------------------
# Timeline sheet: fill in newly finished Actual/Finished dates & notes
# ---------------------------------------------------------------------

# Row 7 - "Define Technical Design Document": Finished 15-Oct, note added
$ws.Range("F7").Value = 43023
$ws.Range("F7").NumberFormat = $dateFmt
$ws.Range("G7").Value = "Switch Priority to project structure"

# Row 8 - "Implementation Phase" header: Actual start date filled in
$ws.Range("E8").Value = 43014
$ws.Range("E8").NumberFormat = $dateFmt

# Row 9 - "Create Solution project scratch": Actual start/finish
$ws.Range("E9").Value = 43014
$ws.Range("E9").NumberFormat = $dateFmt
$ws.Range("F9").Value = 43014
$ws.Range("F9").NumberFormat = $dateFmt

# Row 10 - "Prepare development libraries": Actual start/finish + note
$ws.Range("E10").Value = 43014
$ws.Range("E10").NumberFormat = $dateFmt
$ws.Range("F10").Value = 43023
$ws.Range("F10").NumberFormat = $dateFmt
$ws.Range("G10").Value = "Stuck on troubleshoot using the library"

# Row 11 - "Implement library sample project": Actual start/finish
$ws.Range("E11").Value = 43023
$ws.Range("E11").NumberFormat = $dateFmt
$ws.Range("F11").Value = 43023
$ws.Range("F11").NumberFormat = $dateFmt

# Row 16 - "Implement VPA Face Recognition API Module" / Tests: Actual start + note
$ws.Range("E16").Value = 43023
$ws.Range("E16").NumberFormat = $dateFmt
$ws.Range("G16").Value = "Start Unit Test after development finished (Each module)"

# Update the view so the freshly-edited area is in view
$ws.Range("G17").Select()

# ---------------------------------------------------------------------
# New sheet: Development Module Breakdown
# ---------------------------------------------------------------------

$ws2 = $wb.Worksheets.Add($null, $ws)
$ws2.Name = "Development Module Breakdown"

$ws2.Columns.Item(1).ColumnWidth = 43.85546875
$ws2.Columns.Item(2).ColumnWidth = 16.85546875
$ws2.Columns.Item(3).ColumnWidth = 15.140625
$ws2.Columns.Item(4).ColumnWidth = 16.42578125
$ws2.Columns.Item(5).ColumnWidth = 15.7109375

# Header row
$ws2.Range("A1").Value = "Module Part"
$ws2.Range("B1").Value = "Start"
$ws2.Range("C1").Value = "End"
$ws2.Range("A2").Value = "Encryption"
$ws2.Range("B2").Value = 43023
$ws2.Range("B2").NumberFormat = $dateFmt
$ws2.Range("C2").Value = 43023
$ws2.Range("C2").NumberFormat = $dateFmt
$ws2.Range("A3").Value = "I/O Manager"
$ws2.Range("B3").Value = 43023
$ws2.Range("B3").NumberFormat = $dateFmt
$ws2.Range("A4").Value = "Setting Management"

$ws2.Range("E1").Value = "Status"
$ws2.Range("E2").Value = "Passed"

$ws2.Range("D1").Value = "Unit Test"
$ws2.Range("D2").Value = 43023
$ws2.Range("D2").NumberFormat = $dateFmt

$headerRange = $ws2.Range("A1:E1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4108

$ws2.Range("C3").Select()

# Keep the Timeline sheet as the active tab
$ws.Activate()
